$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the login-steps test data in D2 with the Subha123 credentials
$ws.Range("D2").Value = "1. Launch browser and navigate to https://demo.guru99.com/V4/index.php`n2. Enter the valid user name: Subha123`n3. Enter the valid password: subha@123`n4. Click the login button"

# Update the selected cell shown when the workbook is reopened
$ws.Range("D6").Select()
